# ============================================================
# Commit: feat: add 2022-Q4 data
# - Insert a new worksheet "2022-Q4" (fund holdings detail) right
#   after the "总计" (summary) sheet, pushing the existing quarter
#   sheets ("2021-Q4","2021-Q2","2020-Q4") one slot to the right.
# - Add a corresponding summary row to "总计" for 2022-Q4 with
#   holdings count = 9 and market value = 4.75 (亿元).
#
# NOTE on this COM shim's quirks (discovered empirically):
#  * A Worksheet reference obtained via Worksheets.Item(n) is
#    positional; it must be re-fetched after any call that adds /
#    inserts sheets or rows, otherwise it silently resolves to
#    whatever now sits at that slot.
#  * Worksheets.Add(), Range.Insert()/Rows.Insert(), and
#    Range.ClearFormats() all clear the clipboard, so every
#    Copy() must be immediately followed by its PasteSpecial()
#    with no other such call in between.
# ============================================================

$wb = $excel.ActiveWorkbook
$total = $wb.Worksheets.Item(1)          # "总计" summary sheet

# ------------------------------------------------------------
# Step 1: insert the new "2022-Q4" worksheet right after "总计"
# ------------------------------------------------------------
$new = $wb.Worksheets.Add($null, $total)
$new.Name = "2022-Q4"

# Re-fetch the template sheet ("2021-Q4") fresh - it has shifted to
# index 3 now that the new sheet was inserted at index 2.
$templateQ4 = $wb.Worksheets.Item(3)

# Header row (B1:H1): copy the bold/bordered/centered formatting from the template sheet
$templateQ4.Range("B1:H1").Copy()
$new.Range("B1:H1").PasteSpecial(-4122)   # xlPasteFormats

# Column A (index cells for rows 2-10): same style, copied from the template's A2
$templateQ4.Range("A2").Copy()
$new.Range("A2:A10").PasteSpecial(-4122)  # xlPasteFormats

# Now fill in the header text
$new.Range("B1").Value = "基金代码"
$new.Range("C1").Value = "基金名称"
$new.Range("D1").Value = "基金规模"
$new.Range("E1").Value = "股票总仓位"
$new.Range("F1").Value = "仓位占比"
$new.Range("G1").Value = "持有市值(亿元)"
$new.Range("H1").Value = "仓位排名"

# Data rows 2-10 (fund holdings for 2022-Q4).
# Columns B-G are stored as text in the source data (fund codes keep
# leading zeros, percentages keep trailing zeros), so force a Text
# number format before writing the value, then drop the format again
# so the cell is left unstyled like the rest of the sheet (this also
# matches the source: none of these data cells carry an "s" style).
# row 2
$new.Range("A2").Value = 0
$new.Range("B2").NumberFormat = "@"
$new.Range("B2").Value = "270023"
$new.Range("B2").ClearFormats()
$new.Range("C2").NumberFormat = "@"
$new.Range("C2").Value = "广发全球精选股票（QDII）"
$new.Range("C2").ClearFormats()
$new.Range("D2").NumberFormat = "@"
$new.Range("D2").Value = "20.45"
$new.Range("D2").ClearFormats()
$new.Range("E2").NumberFormat = "@"
$new.Range("E2").Value = "82.63"
$new.Range("E2").ClearFormats()
$new.Range("F2").NumberFormat = "@"
$new.Range("F2").Value = "6.10"
$new.Range("F2").ClearFormats()
$new.Range("G2").NumberFormat = "@"
$new.Range("G2").Value = "1.2474"
$new.Range("G2").ClearFormats()
$new.Range("H2").Value = 6

# row 3
$new.Range("A3").Value = 1
$new.Range("B3").NumberFormat = "@"
$new.Range("B3").Value = "000906"
$new.Range("B3").ClearFormats()
$new.Range("C3").NumberFormat = "@"
$new.Range("C3").Value = "广发全球精选股票（QDII）美元现汇"
$new.Range("C3").ClearFormats()
$new.Range("D3").NumberFormat = "@"
$new.Range("D3").Value = "20.45"
$new.Range("D3").ClearFormats()
$new.Range("E3").NumberFormat = "@"
$new.Range("E3").Value = "82.63"
$new.Range("E3").ClearFormats()
$new.Range("F3").NumberFormat = "@"
$new.Range("F3").Value = "6.10"
$new.Range("F3").ClearFormats()
$new.Range("G3").NumberFormat = "@"
$new.Range("G3").Value = "1.2474"
$new.Range("G3").ClearFormats()
$new.Range("H3").Value = 6

# row 4
$new.Range("A4").Value = 2
$new.Range("B4").NumberFormat = "@"
$new.Range("B4").Value = "011423"
$new.Range("B4").ClearFormats()
$new.Range("C4").NumberFormat = "@"
$new.Range("C4").Value = "广发全球科技三个月定期开放混合（QDII）美元 C"
$new.Range("C4").ClearFormats()
$new.Range("D4").NumberFormat = "@"
$new.Range("D4").Value = "25.66"
$new.Range("D4").ClearFormats()
$new.Range("E4").NumberFormat = "@"
$new.Range("E4").Value = "89.07"
$new.Range("E4").ClearFormats()
$new.Range("F4").NumberFormat = "@"
$new.Range("F4").Value = "4.16"
$new.Range("F4").ClearFormats()
$new.Range("G4").NumberFormat = "@"
$new.Range("G4").Value = "1.0675"
$new.Range("G4").ClearFormats()
$new.Range("H4").Value = 8

# row 5
$new.Range("A5").Value = 3
$new.Range("B5").NumberFormat = "@"
$new.Range("B5").Value = "011420"
$new.Range("B5").ClearFormats()
$new.Range("C5").NumberFormat = "@"
$new.Range("C5").Value = "广发全球科技三个月定期开放混合（QDII）人民币 A"
$new.Range("C5").ClearFormats()
$new.Range("D5").NumberFormat = "@"
$new.Range("D5").Value = "21.02"
$new.Range("D5").ClearFormats()
$new.Range("E5").NumberFormat = "@"
$new.Range("E5").Value = "89.07"
$new.Range("E5").ClearFormats()
$new.Range("F5").NumberFormat = "@"
$new.Range("F5").Value = "4.16"
$new.Range("F5").ClearFormats()
$new.Range("G5").NumberFormat = "@"
$new.Range("G5").Value = "0.8744"
$new.Range("G5").ClearFormats()
$new.Range("H5").Value = 8

# row 6
$new.Range("A6").Value = 4
$new.Range("B6").NumberFormat = "@"
$new.Range("B6").Value = "011422"
$new.Range("B6").ClearFormats()
$new.Range("C6").NumberFormat = "@"
$new.Range("C6").Value = "广发全球科技三个月定期开放混合（QDII）人民币 C"
$new.Range("C6").ClearFormats()
$new.Range("D6").NumberFormat = "@"
$new.Range("D6").Value = "4.84"
$new.Range("D6").ClearFormats()
$new.Range("E6").NumberFormat = "@"
$new.Range("E6").Value = "89.07"
$new.Range("E6").ClearFormats()
$new.Range("F6").NumberFormat = "@"
$new.Range("F6").Value = "4.16"
$new.Range("F6").ClearFormats()
$new.Range("G6").NumberFormat = "@"
$new.Range("G6").Value = "0.2013"
$new.Range("G6").ClearFormats()
$new.Range("H6").Value = 8

# row 7
$new.Range("A7").Value = 5
$new.Range("B7").NumberFormat = "@"
$new.Range("B7").Value = "006792"
$new.Range("B7").ClearFormats()
$new.Range("C7").NumberFormat = "@"
$new.Range("C7").Value = "鹏华香港美国互联网股票（LOF）美元现汇"
$new.Range("C7").ClearFormats()
$new.Range("D7").NumberFormat = "@"
$new.Range("D7").Value = "1.29"
$new.Range("D7").ClearFormats()
$new.Range("E7").NumberFormat = "@"
$new.Range("E7").Value = "88.46"
$new.Range("E7").ClearFormats()
$new.Range("F7").NumberFormat = "@"
$new.Range("F7").Value = "2.80"
$new.Range("F7").ClearFormats()
$new.Range("G7").NumberFormat = "@"
$new.Range("G7").Value = "0.0361"
$new.Range("G7").ClearFormats()
$new.Range("H7").Value = 10

# row 8
$new.Range("A8").Value = 6
$new.Range("B8").NumberFormat = "@"
$new.Range("B8").Value = "160644"
$new.Range("B8").ClearFormats()
$new.Range("C8").NumberFormat = "@"
$new.Range("C8").Value = "鹏华香港美国互联网股票（LOF）人民币"
$new.Range("C8").ClearFormats()
$new.Range("D8").NumberFormat = "@"
$new.Range("D8").Value = "1.29"
$new.Range("D8").ClearFormats()
$new.Range("E8").NumberFormat = "@"
$new.Range("E8").Value = "88.46"
$new.Range("E8").ClearFormats()
$new.Range("F8").NumberFormat = "@"
$new.Range("F8").Value = "2.80"
$new.Range("F8").ClearFormats()
$new.Range("G8").NumberFormat = "@"
$new.Range("G8").Value = "0.0361"
$new.Range("G8").ClearFormats()
$new.Range("H8").Value = 10

# row 9
$new.Range("A9").Value = 7
$new.Range("B9").NumberFormat = "@"
$new.Range("B9").Value = "005698"
$new.Range("B9").ClearFormats()
$new.Range("C9").NumberFormat = "@"
$new.Range("C9").Value = "华夏全球科技先锋混合（QDII）"
$new.Range("C9").ClearFormats()
$new.Range("D9").NumberFormat = "@"
$new.Range("D9").Value = "0.60"
$new.Range("D9").ClearFormats()
$new.Range("E9").NumberFormat = "@"
$new.Range("E9").Value = "83.35"
$new.Range("E9").ClearFormats()
$new.Range("F9").NumberFormat = "@"
$new.Range("F9").Value = "4.96"
$new.Range("F9").ClearFormats()
$new.Range("G9").NumberFormat = "@"
$new.Range("G9").Value = "0.0298"
$new.Range("G9").ClearFormats()
$new.Range("H9").Value = 6

# row 10
$new.Range("A10").Value = 8
$new.Range("B10").NumberFormat = "@"
$new.Range("B10").Value = "011421"
$new.Range("B10").ClearFormats()
$new.Range("C10").NumberFormat = "@"
$new.Range("C10").Value = "广发全球科技三个月定期开放混合（QDII）美元 A"
$new.Range("C10").ClearFormats()
$new.Range("D10").NumberFormat = "@"
$new.Range("D10").Value = "0.20"
$new.Range("D10").ClearFormats()
$new.Range("E10").NumberFormat = "@"
$new.Range("E10").Value = "89.07"
$new.Range("E10").ClearFormats()
$new.Range("F10").NumberFormat = "@"
$new.Range("F10").Value = "4.16"
$new.Range("F10").ClearFormats()
$new.Range("G10").NumberFormat = "@"
$new.Range("G10").Value = "0.0083"
$new.Range("G10").ClearFormats()
$new.Range("H10").Value = 8

# ------------------------------------------------------------
# Step 2: update the "总计" summary sheet - insert a row for
# 2022-Q4 at the top of the data (row 2), shifting the existing
# 2021-Q4 / 2021-Q2 / 2020-Q4 rows down, and renumber column A.
# ------------------------------------------------------------
# Re-fetch fresh: "总计" stays at index 1 (it's before the insert
# point), but re-resolve anyway to be safe since the worksheet
# collection changed shape in Step 1.
$total = $wb.Worksheets.Item(1)

$total.Rows.Item(2).Insert()

# The freshly inserted row inherits the header row's bold format
# on columns B:D; strip that back to the plain (unstyled) look
# used by the other data rows.
$total.Range("B2:D2").ClearFormats()

# Give A2 the same bold/bordered/centered style as the other
# index cells in column A (copy format from A3, which is the old
# A2 pushed down by the insert).
$total.Range("A3").Copy()
$total.Range("A2").PasteSpecial(-4122)   # xlPasteFormats

$total.Range("A2").Value = 0
$total.Range("B2").Value = "2022-Q4"
$total.Range("C2").Value = 9
$total.Range("D2").Value = 4.75

# Renumber the remaining index column (0,1,2,3 top to bottom)
$total.Range("A3").Value = 1
$total.Range("A4").Value = 2
$total.Range("A5").Value = 3
